# Applies the cryptos.xlsx price/volume update described in the commit
# "Updated cryptos list on Tue Jun 25 16:08:29 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.852.39"
$ws.Range("E2").Value = "  +1.86%  "
$ws.Range("D3").Value = "3.420.77"
$ws.Range("E3").Value = "  +4.28%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'577.60"
$ws.Range("E5").Value = "  +2.85%  "
$ws.Range("D6").Value = "'139.49"
$ws.Range("E6").Value = "  +10.34%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.420.85"
$ws.Range("E8").Value = "  +4.31%  "
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("D10").Value = "'7.71"
$ws.Range("E10").Value = "  +6.58%  "
$ws.Range("D11").Value = "'0.127"
$ws.Range("E11").Value = "  +8.60%  "
$ws.Range("E12").Value = "  +6.72%  "
$ws.Range("D13").Value = "3.996.88"
$ws.Range("E13").Value = "  +4.22%  "
$ws.Range("D14").Value = "'0.121"
$ws.Range("E14").Value = "  +1.88%  "
$ws.Range("E15").Value = "  +9.38%  "
$ws.Range("D16").Value = "3.414.46"
$ws.Range("E16").Value = "  +4.11%  "
$ws.Range("D17").Value = "'25.66"
$ws.Range("E17").Value = "  +7.00%  "
$ws.Range("D18").Value = "61.908.68"
$ws.Range("E18").Value = "  +1.66%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "'5.95"
$ws.Range("E19").Value = "  +7.00%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'14.10"
$ws.Range("E20").Value = "  +6.90%  "
$ws.Range("D21").Value = "'9.44"
$ws.Range("E21").Value = "  +6.78%  "
$ws.Range("D22").Value = "'392.17"
$ws.Range("E22").Value = "  +12.20%  "
$ws.Range("E23").Value = "  +4.49%  "
$ws.Range("D24").Value = "3.553.38"
$ws.Range("E24").Value = "  +4.17%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "'0.0000128"
$ws.Range("E25").Value = "  +19.85%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "'71.37"
$ws.Range("E27").Value = "  +3.58%  "
$ws.Range("D28").Value = "'1.64"
$ws.Range("E28").Value = "  +15.61%  "
$ws.Range("D29").Value = "'7.82"
$ws.Range("E29").Value = "  +10.86%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("E31").Value = "  +7.80%  "
$ws.Range("D32").Value = "'0.160"
$ws.Range("E32").Value = "  +8.12%  "
$ws.Range("D33").Value = "'2.17"
$ws.Range("E33").Value = "  +3.77%  "
$ws.Range("D34").Value = "3.446.86"
$ws.Range("E34").Value = "  +4.15%  "
$ws.Range("D36").Value = "'23.69"
$ws.Range("E36").Value = "  +5.31%  "
$ws.Range("D37").Value = "'5.57"
$ws.Range("E37").Value = "  +7.36%  "
$ws.Range("D38").Value = "'7.05"
$ws.Range("E38").Value = "  +4.72%  "
$ws.Range("E39").Value = "  +6.94%  "
$ws.Range("E40").Value = "  +0.84%  "
$ws.Range("E41").Value = "  +7.31%  "
$ws.Range("E42").Value = "  +12.65%  "
$ws.Range("B43").Value = "ONDO"
$ws.Range("C43").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D43").Value = "'1.24"
$ws.Range("E43").Value = "  +11.07%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.777"
$ws.Range("E45").Value = "  +5.45%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "'4.49"
$ws.Range("E46").Value = "  +4.33%  "
$ws.Range("D47").Value = "'41.00"
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("D48").Value = "'23.51"
$ws.Range("E48").Value = "  +7.00%  "
$ws.Range("D49").Value = "'7.01"
$ws.Range("E49").Value = "  +5.44%  "
$ws.Range("D50").Value = "'23.03"
$ws.Range("E50").Value = "  +9.96%  "
$ws.Range("D51").Value = "2.369.63"
$ws.Range("E51").Value = "  +10.64%  "
